$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the text of a paragraph (already fetched as $para) with one
# or more lines of text (an array means several paragraphs replace the one
# original paragraph).
#
# Notes on quirks of this COM host that this helper works around:
#   * Range.Delete() reliably removes text that spans several runs (e.g. a
#     paragraph broken up by <w:proofErr/> spell/grammar markers) as long as
#     the range was not created over text that was *just* inserted together
#     with an embedded carriage return - so we always normalize a paragraph
#     down to a single run first (via a one-line placeholder) and only then
#     perform the (possibly) multi-line insert.
#   * Range.InsertAfter() on a collapsed range that starts at document
#     offset 0 misbehaves (it appends at the end of the story instead of at
#     the collapsed point), so InsertBefore() is used everywhere instead,
#     which behaves correctly at offset 0 too.
# ---------------------------------------------------------------------------
function Replace-ParaText($para, [string[]]$lines) {
    $origStart = $para.Range.Start
    $origEnd = $para.Range.End - 1

    # Step A: normalize the paragraph to a single run containing a unique
    # placeholder, regardless of how many runs/proofErr spans it currently has.
    $placeholder = "PLACEHOLDER-" + [guid]::NewGuid().ToString("N")
    $rIns = $d.Range($origStart, $origStart)
    $rIns.InsertBefore($placeholder)
    $rOld = $d.Range($origStart + $placeholder.Length, $origEnd + $placeholder.Length)
    $rOld.Delete()

    # Step B: replace the placeholder with the real text (one or more lines).
    $newText = [string]::Join("`r", $lines)
    $rIns2 = $d.Range($origStart, $origStart)
    $rIns2.InsertBefore($newText)
    $rOld2 = $d.Range($origStart + $newText.Length, $origStart + $newText.Length + $placeholder.Length)
    $rOld2.Delete()
}

# ---------------------------------------------------------------------------
# Apply edits from the bottom of the document upward so paragraph indices
# for not-yet-processed (earlier) paragraphs stay valid.
# ---------------------------------------------------------------------------

# Paragraph 30 (iframe embed) + paragraph 31 (the blank line right after it)
# collapse into a single new paragraph with the <YouTube .../> import usage.
$pEmptyAfterIframe = $d.Paragraphs(31)
$pEmptyAfterIframe.Range.Delete()
Replace-ParaText $d.Paragraphs(30) @('<YouTube embedId="AJjT5bYknls" />')

# "In 2018, Joseph Briggs ..." - same text, just simplified runs (no proofErr).
Replace-ParaText $d.Paragraphs(28) @('In 2018, Joseph Briggs independently proposed the same concept as an application of the 42 method to the 2x2x2.')

# "![](img/CLL/CCLL.png)" - same text, just simplified runs.
Replace-ParaText $d.Paragraphs(26) @('![](img/CLL/CCLL.png)')

# SpeedSolving wiki link line - same text, just simplified runs.
Replace-ParaText $d.Paragraphs(18) @('[Click here for more step details on the SpeedSolving wiki](https://www.speedsolving.com/wiki/index.php?title=Conjugated_CxLL)')

# "**Steps:**" - same text, just simplified runs.
Replace-ParaText $d.Paragraphs(13) @('**Steps:**')

# "**Proposed:** 2012" - same text, just simplified runs.
Replace-ParaText $d.Paragraphs(11) @('**Proposed:** 2012')

# "**Proposer:** ..." - add the Joseph Briggs credit.
Replace-ParaText $d.Paragraphs(9) @('**Proposer:** [Michael James Straughan](CubingContributors/MethodDevelopers.md#straughan-michael-james-athefre), [Joseph Briggs](CubingContributors/MethodDevelopers.md#briggs-joseph-shadowslice)')

# <AnimCube2x2 .../> viewer embed -> multi-line <ReconViewer ... /> usage.
Replace-ParaText $d.Paragraphs(5) @(
    '<ReconViewer',
    'puzzle="2x2x2"',
    "scramble=""U R F U R U' R' F'""",
    "solution={``F R U R' U' F' . // CCLL",
    "R' U' // Undo pseudo``}",
    '/>'
)

# Front matter / imports at the very top of the file.
Replace-ParaText $d.Paragraphs(1) @(
    '---',
    "description: History of the CCLL method for the 2x2 Rubik's Cube.",
    '---',
    '',
    'import ReconViewer from "@site/src/components/ReconViewer";',
    'import YouTube from "@site/src/components/YouTube";'
)
